$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 4 (tyre id 3) with newly fitted model values ---
$ws.Range("C4").Value = -11.363636363636363
$ws.Range("D4").Value = 1378.5647730407597
$ws.Range("E4").Value = -11.363636363636363
$ws.Range("F4").Value = 1152.3653094664892
$ws.Range("G4").Value = 19143.469473247831
$ws.Range("H4").Value = 16469.351270755647

# --- Fill in previously-empty row 13 (tyre id 12) with fitted / graphed results ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "12_Hoosier_20.5x7-13_R20_8Rim.tir"
$ws.Range("C13").Value = -5.9090909090909092
$ws.Range("D13").Value = 1446.814738053685
$ws.Range("E13").Value = -5.6060606060606055
$ws.Range("F13").Value = 1213.769009670656
$ws.Range("G13").Value = 36968.755564889514
$ws.Range("H13").Value = 31950.493206296178
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1286.6649923608313
$ws.Range("M13").Value = 7.7272727272727275
$ws.Range("N13").Value = 1073.7162945574842
$ws.Range("O13").Value = 8.0303030303030294

$wb.Save()
